$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (cases) values that changed for existing rows
$ws.Range("B21").Value = 2
$ws.Range("B22").Value = 4
$ws.Range("B28").Value = 3
$ws.Range("B29").Value = 4
$ws.Range("B32").Value = 7
$ws.Range("B35").Value = 11
$ws.Range("B37").Value = 8
$ws.Range("B38").Value = 23
$ws.Range("B40").Value = 16
$ws.Range("B42").Value = 36
$ws.Range("B43").Value = 34
$ws.Range("B44").Value = 43
$ws.Range("B45").Value = 86
$ws.Range("B46").Value = 65
$ws.Range("B47").Value = 106
$ws.Range("B48").Value = 100
$ws.Range("B49").Value = 161
$ws.Range("B50").Value = 133
$ws.Range("B51").Value = 319
$ws.Range("B52").Value = 309
$ws.Range("B53").Value = 358
$ws.Range("B54").Value = 399
$ws.Range("B55").Value = 473
$ws.Range("B56").Value = 605
$ws.Range("B57").Value = 687
$ws.Range("B58").Value = 984
$ws.Range("B59").Value = 1599
$ws.Range("B60").Value = 1760
$ws.Range("B61").Value = 2193
$ws.Range("B62").Value = 3122
$ws.Range("B63").Value = 4123
$ws.Range("B64").Value = 4552
$ws.Range("B65").Value = 5554
$ws.Range("B66").Value = 6844
$ws.Range("B67").Value = 6057
$ws.Range("B68").Value = 5886
$ws.Range("B69").Value = 7481
$ws.Range("B70").Value = 8430
$ws.Range("B71").Value = 8440
$ws.Range("B72").Value = 8836
$ws.Range("B73").Value = 9529
$ws.Range("B74").Value = 6854
$ws.Range("B75").Value = 6348
$ws.Range("B76").Value = 9965
$ws.Range("B77").Value = 9788
$ws.Range("B78").Value = 9378
$ws.Range("B79").Value = 8738
$ws.Range("B80").Value = 7950
$ws.Range("B81").Value = 4352
$ws.Range("B82").Value = 2255
$ws.Range("B83").Value = 1185
$ws.Range("B84").Value = 243
$ws.Range("B85").Value = 132
$ws.Range("B86").Value = 47
$ws.Range("B87").Value = 9

# A77 loses its highlight fill (was the first highlighted row; highlight now starts at A78)
$ws.Range("A77").Interior.ColorIndex = -4142
$ws.Range("A77").Interior.Pattern = -4142

# Append new row 88 (2020-04-08 data point)
$ws.Range("A88").Value = 43928
$ws.Range("A88").Interior.Color = 65535
$ws.Range("B88").Value = 0

# Reset selection to the top-left cell
$ws.Range("A1").Select() | Out-Null
